$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("A1").Value = "time"

# New column E header - copy the existing header formatting (bold, borders,
# centered) from D1 onto E1 before setting its text.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "Praat Label"

# Per-row updated values for columns B (Avg Pitch), C (Avg Intensity), and
# new column E (Praat Label). Row 29's B value is cleared (no pitch value).
$rows = @(
    @{Row=2; B=191.6896484019795; C=61.18225459059844; E="joy"}
    @{Row=3; B=231.8624386244848; C=52.53484972441996; E="fear"}
    @{Row=4; B=233.5510455051854; C=65.27336074361563; E="fear"}
    @{Row=5; B=177.2476450762052; C=65.77709160669976; E="fear"}
    @{Row=6; B=185.882137541738; C=57.21560925996133; E="joy"}
    @{Row=7; B=216.4374150207751; C=61.3490178925293; E="fear"}
    @{Row=8; B=167.4480746313774; C=54.70533428075144; E="joy"}
    @{Row=9; B=201.8726402485668; C=66.29160810830206; E="fear"}
    @{Row=10; B=245.6079890769943; C=53.22962207934772; E="fear"}
    @{Row=11; B=233.823800445258; C=53.27962219302967; E="fear"}
    @{Row=12; B=247.3240934338943; C=65.58672811893763; E="fear"}
    @{Row=13; B=208.1576955976502; C=58.33563423240206; E="fear"}
    @{Row=14; B=201.3080359845419; C=68.31256952537977; E="fear"}
    @{Row=15; B=248.6866536445313; C=65.39013064218855; E="fear"}
    @{Row=16; B=246.8835872138402; C=67.6579554849902; E="fear"}
    @{Row=17; B=204.4413285475489; C=63.95552118291549; E="fear"}
    @{Row=18; B=236.9850972149102; C=66.84129124089847; E="fear"}
    @{Row=19; B=252.4475508835506; C=60.94675076648562; E="fear"}
    @{Row=20; B=227.1528452591361; C=70.48025136775708; E="fear"}
    @{Row=21; B=180.1942825438341; C=63.82513395829859; E="fear"}
    @{Row=22; B=227.5726412942738; C=70.85039711827095; E="fear"}
    @{Row=23; B=239.7605994135779; C=68.40482233510463; E="fear"}
    @{Row=24; B=255.4938470504626; C=69.31847242164309; E="fear"}
    @{Row=25; B=220.8254725672719; C=70.437806301094; E="fear"}
    @{Row=26; B=202.5392092414347; C=70.90531984436423; E="fear"}
    @{Row=27; B=231.7528436958735; C=64.27213797476745; E="fear"}
    @{Row=28; B=248.5583206852246; C=65.56225596180012; E="fear"}
    @{Row=29; B=$null; C=48.72176356844336; E="sadness"}
    @{Row=30; B=155.3043119940262; C=53.09575079729106; E="joy"}
    @{Row=31; B=225.1850319179004; C=57.61391680600597; E="fear"}
)

foreach ($r in $rows) {
    if ($null -eq $r.B) {
        $ws.Cells.Item($r.Row, 2).Value = $null
    } else {
        $ws.Cells.Item($r.Row, 2).Value = $r.B
    }
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
